# Ratios.xlsx - corrected age calculation algorithm:
#  - recompute dU234 (col A) and its absolute error (col B) for rows 2-12
#  - column B was narrowed slightly to fit the new formatting

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column B (dU234 error column) from 23.7109375 to 22.7109375 characters.
# (The host's ColumnWidth setter quantizes to the nearest 1/6-character grid point
# before writing the OOXML <col width>, so 21.833333333333332 is the input that lands
# closest to the target 22.7109375 once that fixed +5/6 offset is applied.)
$ws.Columns.Item(2).ColumnWidth = 21.833333333333332

# Recomputed values: column A = dU234, column B = Error dU234 (abs.)
$values = @{
    2  = @(-1.385999438767938,  0.001669580863432305)
    3  = @(144.2958033305279,   0.0005475608978145764)
    4  = @(-1.425063738750176,  0.001844886627987924)
    5  = @(144.2374854732194,   0.0004965169137167292)
    6  = @(1.010673376144799,   0.001799454776365855)
    7  = @(144.6389429527974,   0.0005465452420743875)
    8  = @(-0.8032792008549894, 0.001331396634869054)
    9  = @(144.2738303060913,   0.0003788097800316013)
    10 = @(-0.7966043493090602, 0.001926555786283)
    11 = @(145.2398038022811,   0.0004874012285628684)
    12 = @(-0.3808579399658951, 0.00150393304762386)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
}
